$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Insert new row "FA22 Reindizierung anstoßen" right before "NFA01 Performanz Suche" ---
$targetRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13) -eq "NFA01 Performanz Suche") {
        $targetRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($targetRow)
$newIndex = $newRow.Index
$t.Cell($newIndex, 1).Range.Text = "FA22 Reindizierung anstoßen"
$t.Cell($newIndex, 2).Range.Text = "Über die Anbindung des Webinterfaces der Suchengine kann der Admin eine Reindizierung manuell anstoßen. Anschließend führt das System die Reindizierung durch."
$t.Cell($newIndex, 3).Range.Text = "JA"

# --- 2. Move the <w:lastRenderedPageBreak/> marker from the "NFA04" row to the "NFA02" row ---
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$nfa02Row = $null
$nfa04Row = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $txt = $t.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13)
    if ($txt -eq "NFA02 Performanz Reindizieren") {
        $nfa02Row = $i
    }
    if ($txt -eq "NFA04 Bedienbarkeit Einarbeitungszeit") {
        $nfa04Row = $i
    }
}

# remove it from the NFA04 cell
$nfa04Cell = $t.Cell($nfa04Row, 1)
$nfa04Cell.Range.InsertXML("<w:p $wNs><w:r><w:t>NFA04 Bedienbarkeit Einarbeitungszeit</w:t></w:r></w:p>")

# add it to the NFA02 cell
$nfa02Cell = $t.Cell($nfa02Row, 1)
$nfa02Cell.Range.InsertXML("<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>NFA02 Performanz Reindizieren</w:t></w:r></w:p>")
